$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" text on Hoja1 (A1) ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.32 = 8787.91 pesos`n✅ 8787.91 pesos = 2.3 = 958.92 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# --- Update the rate figures on the "tasas" sheet ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 430.99
$ws2.Range("O10").Value = 3787.5

$ws2.Range("N12").Value = 3818
$ws2.Range("O12").Value = 416.611
